# Update registered/benefit-recipient family statistics on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("registered family") - update years 2015-2021 (columns E-K)
$ws.Range("E4").Value = 4182
$ws.Range("F4").Value = 3652
$ws.Range("G4").Value = 3446
$ws.Range("H4").Value = 3424
$ws.Range("I4").Value = 3358
$ws.Range("J4").Value = 3434
$ws.Range("K4").Value = 3471

# Row 5 ("subsistence allowance recipient family") - update years 2015-2021 (columns E-K)
$ws.Range("E5").Value = 1937
$ws.Range("F5").Value = 1956
$ws.Range("G5").Value = 1820
$ws.Range("H5").Value = 1964
$ws.Range("I5").Value = 1701
$ws.Range("J5").Value = 2030
$ws.Range("K5").Value = 2110
